# Update scripts with new TPM data.
# A new "ECs" sending-cluster row is inserted as the new row 2, pushing the
# previously existing Inflammatory-Mac / MuSCs / Resolving-Mac rows down by
# one row, and every row's derived-specificity values are recalculated
# against the new (4-row) totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 so the existing three data rows
# shift down to rows 3, 4 and 5. The insert copies the header row's bold
# bordered formatting onto the new row, so strip that back to the plain
# (unstyled) look used by the other data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# --- Row 2: ECs -> Qrfp -> Qrfpr -> FAPs -------------------------------
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Qrfp"
$ws.Range("C2").Value = "Qrfpr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.155514
$ws.Range("H2").Value = 0.466542
$ws.Range("I2").Value = 0.1307277193751301
$ws.Range("J2").Value = 0.1307277193751301
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01702166666666667
$ws.Range("N2").Value = 0.051065
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00264710747
$ws.Range("R2").Value = 0.02382396723
$ws.Range("S2").Value = 0.1307277193751301
$ws.Range("T2").Value = 0.1307277193751301

# --- Row 3: Inflammatory-Mac -> Qrfp -> Qrfpr -> FAPs ------------------
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Qrfp"
$ws.Range("C3").Value = "Qrfpr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08079
$ws.Range("H3").Value = 0.24237
$ws.Range("I3").Value = 0.06791345118971129
$ws.Range("J3").Value = 0.0679134511897113
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01702166666666667
$ws.Range("N3").Value = 0.051065
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.00137518045
$ws.Range("R3").Value = 0.01237662405
$ws.Range("S3").Value = 0.06791345118971129
$ws.Range("T3").Value = 0.0679134511897113

# --- Row 4: MuSCs -> Qrfp -> Qrfpr -> FAPs ------------------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Qrfp"
$ws.Range("C4").Value = "Qrfpr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5892876666666667
$ws.Range("H4").Value = 1.767863
$ws.Range("I4").Value = 0.4953652579139191
$ws.Range("J4").Value = 0.4953652579139191
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01702166666666667
$ws.Range("N4").Value = 0.051065
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.01003065823277778
$ws.Range("R4").Value = 0.090275924095
$ws.Range("S4").Value = 0.4953652579139191
$ws.Range("T4").Value = 0.4953652579139191

# --- Row 5: Resolving-Mac -> Qrfp -> Qrfpr -> FAPs ----------------------
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Qrfp"
$ws.Range("C5").Value = "Qrfpr"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3640106666666667
$ws.Range("H5").Value = 1.092032
$ws.Range("I5").Value = 0.3059935715212395
$ws.Range("J5").Value = 0.3059935715212395
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01702166666666667
$ws.Range("N5").Value = 0.051065
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.006196068231111112
$ws.Range("R5").Value = 0.05576461408
$ws.Range("S5").Value = 0.3059935715212395
$ws.Range("T5").Value = 0.3059935715212395
